# Apply updated crypto price/volume data (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.602.94"
$ws.Range("E2").Value = "  -2.65%  "

$ws.Range("D3").Value = "1.982.32"
$ws.Range("E3").Value = "  -3.94%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'242.49"
$ws.Range("E5").Value = "  +0.62%  "

$ws.Range("E6").Value = "  -5.60%  "

$ws.Range("D7").Value = "'57.03"
$ws.Range("E7").Value = "  +8.91%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").Value = "'59.44"
$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("E11").Value = "  -2.82%  "

$ws.Range("D12").Value = "'0.103"
$ws.Range("E12").Value = "  -4.24%  "

$ws.Range("D13").Value = "'0.920"
$ws.Range("E13").Value = "  +3.93%  "

$ws.Range("E14").Value = "  -2.94%  "

$ws.Range("D15").Value = "2.273.56"

$ws.Range("D16").Value = "'5.22"
$ws.Range("E16").Value = "  -3.04%  "

$ws.Range("D17").Value = "1.990.51"
$ws.Range("E17").Value = "  -3.80%  "

$ws.Range("D18").Value = "'17.20"
$ws.Range("E18").Value = "  +5.64%  "

$ws.Range("D19").Value = "35.506.47"
$ws.Range("E19").Value = "  -2.67%  "

$ws.Range("D20").Value = "'70.73"
$ws.Range("E20").Value = "  -1.24%  "

$ws.Range("E21").Value = "  -2.72%  "

$ws.Range("D22").Value = "'233.46"
$ws.Range("E22").Value = "  -0.79%  "

$ws.Range("E23").Value = "  -3.78%  "

$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.28"
$ws.Range("E25").Value = "  -3.34%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.33"
$ws.Range("E26").Value = "  +9.90%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'9.12"
$ws.Range("E27").Value = "  -1.11%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'163.40"
$ws.Range("E28").Value = "  +0.22%  "

$ws.Range("D29").Value = "'19.40"
$ws.Range("E29").Value = "  -4.22%  "

$ws.Range("D30").Value = "'0.119"
$ws.Range("E30").Value = "  -3.63%  "

$ws.Range("E31").Value = "  -0.93%  "

$ws.Range("E32").Value = "  -5.04%  "

$ws.Range("E33").Value = "  -0.85%  "

$ws.Range("E34").Value = "  +10.83%  "

$ws.Range("D35").Value = "'4.27"
$ws.Range("E35").Value = "  -6.17%  "

$ws.Range("D36").Value = "'2.37"
$ws.Range("E36").Value = "  +4.78%  "

$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("E38").Value = "  -1.57%  "

$ws.Range("D39").Value = "'4.93"
$ws.Range("E39").Value = "  +1.69%  "

$ws.Range("E40").Value = "  -4.59%  "

$ws.Range("D41").Value = "'2.83"
$ws.Range("E41").Value = "  -1.96%  "

$ws.Range("E42").Value = "  -2.35%  "

$ws.Range("E43").Value = "  -3.77%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'90.94"
$ws.Range("E44").Value = "  -2.89%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0886"
$ws.Range("E45").Value = "  -4.38%  "

$ws.Range("D46").Value = "1.373.89"
$ws.Range("E46").Value = "  -0.89%  "

$ws.Range("D47").Value = "'7.47"
$ws.Range("E47").Value = "  +3.08%  "

$ws.Range("D48").Value = "'15.52"
$ws.Range("E48").Value = "  -0.15%  "

$ws.Range("E49").Value = "  +1.16%  "

$ws.Range("E50").Value = "  -2.59%  "

$ws.Range("D51").Value = "'45.63"
$ws.Range("E51").Value = "  +2.75%  "
